$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 458 - this shifts the existing row 458
# (and everything below it, down through the old last row 498) down by
# one row, growing the used range from A1:R498 to A1:R499.
$ws.Rows.Item(458).Insert()

# Populate the newly inserted row 458 with the new weekly data point.
$ws.Cells.Item(458, 1).Value = 6
$ws.Cells.Item(458, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(458, 3).Value = "Metropolitana"
$ws.Cells.Item(458, 4).Value = 44918
$ws.Cells.Item(458, 5).Value = 13
$ws.Cells.Item(458, 6).Value = 100112032
$ws.Cells.Item(458, 7).Value = "Zapallo italiano"
$ws.Cells.Item(458, 8).Value = "Sin especificar"
$ws.Cells.Item(458, 9).Value = "Primera"
$ws.Cells.Item(458, 10).Value = 400
$ws.Cells.Item(458, 11).Value = 4000
$ws.Cells.Item(458, 12).Value = 5000
$ws.Cells.Item(458, 13).Value = 4425
$ws.Cells.Item(458, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(458, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(458, 16).Value = 88
$ws.Cells.Item(458, 17).Value = 50
$ws.Cells.Item(458, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D (s="2").
$ws.Range("D458").NumberFormat = $ws.Range("D459").NumberFormat
